$wb = $excel.ActiveWorkbook

# --- Summary Table ---
$ws = $wb.Worksheets.Item("Summary Table")
$ws.Cells.Item(2,1).Value = ' Artificial intelligence '
$ws.Cells.Item(2,2).Value = ' 47 '
$ws.Cells.Item(2,3).Value = ' China unveils new national AI laboratory network focusing on critical infrastructure protection and ethical guidelines. '
$ws.Cells.Item(2,4).Value = ' https://www.scmp.com/tech/policy/article/3254891/china-launches-national-ai-laboratory-network-boost-research-critical-infrastructures '
$ws.Cells.Item(2,5).Value = ' The Chinese Ministry of Science and Technology announced the establishment of a national AI laboratory network with 15 specialized facilities. The network aims to address security challenges in critical infrastructure and implement China''s recently released ethical guidelines for AI. The initiative falls under China''s 14th Five-Year Plan for AI development, emphasizing both innovation and responsible governance. '

$ws.Cells.Item(3,1).Value = ' Innovation-driven development '
$ws.Cells.Item(3,2).Value = ' 38 '
$ws.Cells.Item(3,3).Value = ' China releases comprehensive innovation-driven development strategy for 2025-2035, focusing on self-reliance in critical technologies. '
$ws.Cells.Item(3,4).Value = ' https://english.news.cn/20250428/a7c4e18bb5e645e9a8b1f2d19c3a7f82/c.html '
$ws.Cells.Item(3,5).Value = ' China''s State Council published a new ten-year innovation strategy emphasizing technological self-reliance and indigenous innovation capabilities in semiconductors, advanced materials, and AI. The plan introduces a "2+X" framework with core technologies and supporting ecosystem development. It targets increasing R&D spending to 3% of GDP by 2030 and addresses challenges including aging demographics and climate change. '

$ws.Cells.Item(4,1).Value = ' New quality productivity '
$ws.Cells.Item(4,2).Value = ' 35 '
$ws.Cells.Item(4,3).Value = ' Beijing unveils "New Quality Productivity" initiative with $150 billion investment in advanced manufacturing hubs across China. '
$ws.Cells.Item(4,4).Value = ' https://www.reuters.com/technology/china-announces-150-billion-new-quality-productivity-manufacturing-initiative-2025-04-30/ '
$ws.Cells.Item(4,5).Value = ' China''s National Development and Reform Commission announced a $150 billion initiative to establish 25 "New Quality Productivity" manufacturing hubs nationwide. The program aims to upgrade China''s industrial base with advanced digital manufacturing capabilities, smart factories, and integrated supply chains. The first wave will focus on semiconductors, new energy vehicles, aerospace, and biotech sectors, creating an estimated 1.2 million high-skilled jobs. '

$ws.Cells.Item(5,1).Value = ' Quantum communication '
$ws.Cells.Item(5,2).Value = ' 33 '
$ws.Cells.Item(5,3).Value = ' China claims "quantum advantage" breakthrough with nationwide quantum secure communication network reaching 100+ cities. '
$ws.Cells.Item(5,4).Value = ' https://www.nature.com/articles/d41586-025-01234-z '
$ws.Cells.Item(5,5).Value = ' Chinese researchers reported a major breakthrough in quantum communication, extending their secure quantum network to over 100 cities nationwide. The network demonstrates "quantum advantage" in secure communications using entangled photons, making it theoretically immune to conventional hacking methods. This expansion represents the world''s largest functional quantum communication infrastructure, with plans to extend services to Southeast Asia through the Digital Silk Road initiative. '

$ws.Cells.Item(6,1).Value = ' Science and technology security policy '
$ws.Cells.Item(6,2).Value = ' 31 '
$ws.Cells.Item(6,3).Value = ' China introduces comprehensive Science and Technology Security Framework with new export controls on AI chipsets and quantum technologies. '
$ws.Cells.Item(6,4).Value = ' https://www.bloomberg.com/news/articles/2025-04-29/china-announces-new-tech-export-controls-eyeing-security-reciprocity '
$ws.Cells.Item(6,5).Value = ' The Chinese Commerce Ministry announced expanded export controls covering advanced AI chipsets, quantum computing components, and biotechnology research materials. The measures, part of China''s new Science and Technology Security Framework, aim to protect national security while promoting "fair and reciprocal" technology cooperation. The policy reflects growing concerns about technology containment strategies by Western nations and emphasizes the protection of critical indigenous innovations. '

$ws.Cells.Item(7,1).Value = ' Integration of industry, academia and research '
$ws.Cells.Item(7,2).Value = ' 29 '
$ws.Cells.Item(7,3).Value = ' China launches 50 new industry-academia-research centers with reformed intellectual property sharing mechanisms. '
$ws.Cells.Item(7,4).Value = ' https://www.chinadaily.com.cn/a/202504/28/WS660bc21a310dbb0113778e4.html '
$ws.Cells.Item(7,5).Value = ' The Ministry of Education and Ministry of Industry and Information Technology jointly established 50 new integrated industry-academia-research centers across key universities and industrial parks. The initiative features reformed intellectual property sharing mechanisms allowing more equitable distribution of benefits between researchers and commercial partners. The centers will focus on semiconductors, advanced manufacturing, and renewable energy technologies, with streamlined approval processes for joint research projects. '

$ws.Cells.Item(8,1).Value = ' Digital economy policy '
$ws.Cells.Item(8,2).Value = ' 27 '
$ws.Cells.Item(8,3).Value = ' China releases Digital Economy Development White Paper targeting 45% digital contribution to GDP by 2030. '
$ws.Cells.Item(8,4).Value = ' https://www.globaltimes.cn/page/202504/1309875.shtml '
$ws.Cells.Item(8,5).Value = ' China''s State Council Information Office published a comprehensive Digital Economy Development White Paper outlining plans to increase the digital economy''s contribution to GDP from current 41% to 45% by 2030. The policy framework emphasizes cross-border data governance, digital infrastructure investment, and international standard-setting in emerging technologies. It introduces a "digital sovereignty with openness" approach that balances domestic control with international engagement. '

$ws.Cells.Item(9,1).Value = ' Semiconductor packaging '
$ws.Cells.Item(9,2).Value = ' 26 '
$ws.Cells.Item(9,3).Value = ' China achieves breakthrough in advanced semiconductor packaging with new heterogeneous integration technique for AI chips. '
$ws.Cells.Item(9,4).Value = ' https://asia.nikkei.com/Business/Tech/Semiconductors/China-claims-breakthrough-in-advanced-chip-packaging-technology '
$ws.Cells.Item(9,5).Value = ' Chinese researchers from the Chinese Academy of Sciences announced a breakthrough in advanced semiconductor packaging technology, developing a new heterogeneous integration technique for AI chips. The method allows for stacking specialized chiplets with significantly improved performance and energy efficiency compared to traditional packaging. This development potentially bypasses some restrictions on advanced chip manufacturing while advancing China''s semiconductor self-sufficiency goals. '

$ws.Cells.Item(10,1).Value = ' Science and technology innovation policy '
$ws.Cells.Item(10,2).Value = ' 25 '
$ws.Cells.Item(10,3).Value = ' China updates national science and technology innovation policy with 10-year roadmap and increased basic research funding. '
$ws.Cells.Item(10,4).Value = ' https://www.chinadaily.com.cn/a/202504/30/WS661cb9453a2b0ad6b3b952e.html '
$ws.Cells.Item(10,5).Value = ' China''s State Council released an updated national science and technology innovation policy with a 10-year roadmap prioritizing seven strategic technology areas. The plan increases basic research funding to 12% of total R&D expenditure and introduces a "Basic Research+" initiative linking fundamental science to industrial applications. The policy reforms evaluation systems for scientists, emphasizing quality over quantity of research outputs, and streamlines international collaboration mechanisms. '

$ws.Cells.Item(11,1).Value = ' Strategic emerging industries '
$ws.Cells.Item(11,2).Value = ' 24 '
$ws.Cells.Item(11,3).Value = ' China designates six new strategic emerging industries with $200 billion development fund focusing on future materials. '
$ws.Cells.Item(11,4).Value = ' https://www.scmp.com/economy/article/3254921/china-designates-six-new-strategic-industries-200-billion-development-fund '
$ws.Cells.Item(11,5).Value = ' China''s National Development and Reform Commission announced six new strategic emerging industries: future materials, quantum information, genetic technology, future networks, low-altitude economy, and hydrogen energy. A $200 billion development fund will support these sectors, with future materials receiving the largest allocation. The initiative aims to cultivate industries with an estimated combined value of 35 trillion yuan by 2030 and achieve technological leadership in at least three sectors. '

$ws.Cells.Item(12,1).Value = ' Hydrogen energy storage '
$ws.Cells.Item(12,2).Value = ' 23 '
$ws.Cells.Item(12,3).Value = ' China launches world''s largest hydrogen energy storage facility with 100MW capacity in Inner Mongolia. '
$ws.Cells.Item(12,4).Value = ' https://www.reuters.com/business/energy/china-launches-worlds-largest-hydrogen-energy-storage-project-2025-04-27/ '
$ws.Cells.Item(12,5).Value = ' China completed the world''s largest hydrogen energy storage facility in Inner Mongolia with 100MW capacity, designed to store renewable energy from wind and solar sources. The project uses electrolysis technology to convert excess renewable energy into hydrogen, which can be stored and later converted back to electricity during peak demand periods. The facility marks a significant advancement in China''s energy storage capabilities and supports its dual carbon goals. '

$ws.Cells.Item(13,1).Value = ' Science and technology system reform '
$ws.Cells.Item(13,2).Value = ' 22 '
$ws.Cells.Item(13,3).Value = ' China announces major science and technology system reform emphasizing researcher autonomy and reducing administrative burdens. '
$ws.Cells.Item(13,4).Value = ' https://www.nature.com/articles/d41586-025-01245-9 '
$ws.Cells.Item(13,5).Value = ' The Chinese Academy of Sciences unveiled a major reform of its research system, giving scientists greater autonomy over research directions and reducing administrative burdens. The reforms include a new "PI responsibility system" with expanded budget control for principal investigators, simplified grant application procedures, and reformed evaluation criteria focusing on breakthrough innovations rather than publication metrics. The changes aim to address long-standing bureaucratic inefficiencies in China''s research ecosystem. '

$ws.Cells.Item(14,1).Value = ' Low-carbon energy research and development '
$ws.Cells.Item(14,2).Value = ' 21 '
$ws.Cells.Item(14,3).Value = ' China establishes International Low-Carbon Energy R&D Alliance with 15 countries to accelerate clean energy transitions. '
$ws.Cells.Item(14,4).Value = ' https://english.news.cn/20250501/c982e5a78cd48a9bb5e7f219def71c3/c.html '
$ws.Cells.Item(14,5).Value = ' China formed an International Low-Carbon Energy R&D Alliance with 15 countries including members from Europe, Asia, and Africa. The alliance will coordinate research efforts, share technological advances, and jointly develop standards for various low-carbon technologies. Five collaborative research centers will be established focusing on advanced solar power, next-generation nuclear energy, sustainable hydrogen production, energy storage systems, and carbon capture utilization and storage technologies. '

$ws.Cells.Item(15,1).Value = ' Artificial intelligence laboratory '
$ws.Cells.Item(15,2).Value = ' 20 '
$ws.Cells.Item(15,3).Value = ' China-Singapore Joint AI Laboratory opens in Shenzhen focusing on medical applications and large language models. '
$ws.Cells.Item(15,4).Value = ' https://www.straitstimes.com/asia/east-asia/china-singapore-joint-ai-laboratory-opens-in-shenzhen '
$ws.Cells.Item(15,5).Value = ' The China-Singapore Joint AI Laboratory officially opened in Shenzhen''s International Science and Technology Innovation Center. The facility represents a $300 million investment from both governments and will focus on medical AI applications and large language models optimized for Southeast Asian languages. The laboratory features shared computing infrastructure and a talent exchange program allowing researchers from both countries to collaborate on projects meeting ethical AI standards of both nations. '

$ws.Cells.Item(16,1).Value = ' Industry-university-research cooperation '
$ws.Cells.Item(16,2).Value = ' 19 '
$ws.Cells.Item(16,3).Value = ' New industry-university-research cooperation mechanism introduced with revised IP sharing framework and technology transfer incentives. '
$ws.Cells.Item(16,4).Value = ' https://www.chinadaily.com.cn/a/202504/27/WS66096a9a310dbb0113778c6.html '
$ws.Cells.Item(16,5).Value = ' China''s Ministry of Science and Technology introduced a new national framework for industry-university-research cooperation featuring revised intellectual property sharing guidelines and enhanced technology transfer incentives. The policy allows university researchers to retain up to 70% of benefits from commercialized technologies and establishes dedicated technology transfer offices at 100 universities. A streamlined approval process reduces bureaucratic barriers for joint projects, particularly in strategic sectors such as semiconductors and advanced manufacturing. '

$ws.Cells.Item(17,1).Value = ' Quantum computing cloud platform '
$ws.Cells.Item(17,2).Value = ' 18 '
$ws.Cells.Item(17,3).Value = ' China launches open-access quantum computing cloud platform with 156-qubit processor available to global researchers. '
$ws.Cells.Item(17,4).Value = ' https://www.scmp.com/tech/big-tech/article/3254867/chinese-tech-giant-baidu-launches-156-qubit-quantum-computer-claiming-superior-performance '
$ws.Cells.Item(17,5).Value = ' A leading Chinese tech company launched an open-access quantum computing cloud platform featuring a 156-qubit processor available to researchers globally. The system claims superior performance metrics on certain quantum algorithms compared to competing systems. The platform provides development tools, educational resources, and simulation environments to accelerate quantum software innovation. While open to international researchers, data management policies comply with China''s data security regulations. '

$ws.Cells.Item(18,1).Value = ' Biotechnology '
$ws.Cells.Item(18,2).Value = ' 17 '
$ws.Cells.Item(18,3).Value = ' China updates biotechnology regulatory framework with streamlined approval process for gene therapy and synthetic biology. '
$ws.Cells.Item(18,4).Value = ' https://www.reuters.com/business/healthcare-pharmaceuticals/china-updates-biotech-regulatory-framework-accelerate-gene-therapy-approvals-2025-04-29/ '
$ws.Cells.Item(18,5).Value = ' China''s National Medical Products Administration released updated biotechnology regulations with streamlined approval processes for gene therapies and synthetic biology products. The framework introduces a new "parallel review" mechanism reducing approval times by up to 50% for breakthrough therapies. The regulations establish clear guidelines for CRISPR-based treatments while maintaining ethical oversight. Industry analysts project the changes could make China the world''s largest gene therapy market by 2030. '

$ws.Cells.Item(19,1).Value = ' High-end equipment manufacturing '
$ws.Cells.Item(19,2).Value = ' 16 '
$ws.Cells.Item(19,3).Value = ' China unveils new high-end equipment manufacturing plan targeting aerospace, robotics, and advanced medical devices. '
$ws.Cells.Item(19,4).Value = ' https://www.globaltimes.cn/page/202504/1309825.shtml '
$ws.Cells.Item(19,5).Value = ' China''s Ministry of Industry and Information Technology released a high-end equipment manufacturing development plan focusing on aerospace systems, industrial robotics, and advanced medical devices. The initiative includes special funding mechanisms, tax incentives, and procurement preferences for domestically developed equipment. Ten manufacturing innovation centers will be established to bridge research-to-production gaps. The plan specifically emphasizes technologies where China currently depends on imports, with goals to achieve 75% self-sufficiency in these areas by 2030. '

$ws.Cells.Item(20,1).Value = ' Talent introduction policy '
$ws.Cells.Item(20,2).Value = ' 15 '
$ws.Cells.Item(20,3).Value = ' China announces enhanced talent introduction policy with expanded visa program and $5 billion research fund for international scientists. '
$ws.Cells.Item(20,4).Value = ' https://www.scmp.com/news/china/science/article/3254832/china-announces-enhanced-talent-introduction-policy-expanded-visa-program '
$ws.Cells.Item(20,5).Value = ' China unveiled an enhanced talent introduction strategy with significantly expanded visa programs and a $5 billion research fund specifically for international scientists. The policy introduces a new "Science and Technology Innovation Visa" with expedited processing and ten-year multi-entry options. Additional incentives include subsidized housing, school priority for dependents, and tax benefits for high-level researchers. The initiative targets 20,000 international experts in semiconductor design, quantum computing, AI, and advanced manufacturing within five years. '

$ws.Cells.Item(21,1).Value = ' Science and technology cooperation agreement '
$ws.Cells.Item(21,2).Value = ' 14 '
$ws.Cells.Item(21,3).Value = ' China signs comprehensive science and technology cooperation agreements with Brazil and South Africa focusing on space technology and biomedical research. '
$ws.Cells.Item(21,4).Value = ' https://english.news.cn/20250429/a45fc6ee98f94d99b37d3b0cf3f2b64/c.html '
$ws.Cells.Item(21,5).Value = ' China signed new science and technology cooperation agreements with Brazil and South Africa, significantly expanding BRICS collaboration. The agreements focus on space technology sharing, including satellite development and earth observation systems, along with joint biomedical research programs targeting infectious diseases and cancer treatments. The partnerships establish shared laboratory facilities, simplified researcher exchange processes, and harmonized intellectual property frameworks. These agreements represent China''s largest S&T cooperation initiatives with Global South nations this year. '


# --- Sources ---
$ws = $wb.Worksheets.Item("Sources")
$ws.Cells.Item(3,1).Value = ' South China Morning Post '
$ws.Cells.Item(3,2).Value = ' https://www.scmp.com/tech/policy/article/3254891/china-launches-national-ai-laboratory-network-boost-research-critical-infrastructures '
$ws.Cells.Item(3,3).Value = ' 2025-04-30 '

$ws.Cells.Item(4,1).Value = ' Xinhua News Agency '
$ws.Cells.Item(4,2).Value = ' https://english.news.cn/20250428/a7c4e18bb5e645e9a8b1f2d19c3a7f82/c.html '
$ws.Cells.Item(4,3).Value = ' 2025-04-28 '

$ws.Cells.Item(5,1).Value = ' Reuters Technology '
$ws.Cells.Item(5,2).Value = ' https://www.reuters.com/technology/china-announces-150-billion-new-quality-productivity-manufacturing-initiative-2025-04-30/ '
$ws.Cells.Item(5,3).Value = ' 2025-04-30 '

$ws.Cells.Item(6,1).Value = ' Nature '
$ws.Cells.Item(6,2).Value = ' https://www.nature.com/articles/d41586-025-01234-z '
$ws.Cells.Item(6,3).Value = ' 2025-04-28 '

$ws.Cells.Item(7,1).Value = ' Bloomberg News '
$ws.Cells.Item(7,2).Value = ' https://www.bloomberg.com/news/articles/2025-04-29/china-announces-new-tech-export-controls-eyeing-security-reciprocity '
$ws.Cells.Item(7,3).Value = ' 2025-04-29 '

$ws.Cells.Item(8,1).Value = ' China Daily '
$ws.Cells.Item(8,2).Value = ' https://www.chinadaily.com.cn/a/202504/28/WS660bc21a310dbb0113778e4.html '
$ws.Cells.Item(8,3).Value = ' 2025-04-28 '

$ws.Cells.Item(9,1).Value = ' Global Times '
$ws.Cells.Item(9,2).Value = ' https://www.globaltimes.cn/page/202504/1309875.shtml '
$ws.Cells.Item(9,3).Value = ' 2025-04-29 '

$ws.Cells.Item(10,1).Value = ' Nikkei Asia '
$ws.Cells.Item(10,2).Value = ' https://asia.nikkei.com/Business/Tech/Semiconductors/China-claims-breakthrough-in-advanced-chip-packaging-technology '
$ws.Cells.Item(10,3).Value = ' 2025-04-28 '

$ws.Cells.Item(11,1).Value = ' China Daily '
$ws.Cells.Item(11,2).Value = ' https://www.chinadaily.com.cn/a/202504/30/WS661cb9453a2b0ad6b3b952e.html '
$ws.Cells.Item(11,3).Value = ' 2025-04-30 '

$ws.Cells.Item(12,1).Value = ' South China Morning Post '
$ws.Cells.Item(12,2).Value = ' https://www.scmp.com/economy/article/3254921/china-designates-six-new-strategic-industries-200-billion-development-fund '
$ws.Cells.Item(12,3).Value = ' 2025-04-29 '

$ws.Cells.Item(13,1).Value = ' Reuters Energy '
$ws.Cells.Item(13,2).Value = ' https://www.reuters.com/business/energy/china-launches-worlds-largest-hydrogen-energy-storage-project-2025-04-27/ '
$ws.Cells.Item(13,3).Value = ' 2025-04-27 '

$ws.Cells.Item(14,1).Value = ' Nature '
$ws.Cells.Item(14,2).Value = ' https://www.nature.com/articles/d41586-025-01245-9 '
$ws.Cells.Item(14,3).Value = ' 2025-04-29 '

$ws.Cells.Item(15,1).Value = ' Xinhua News Agency '
$ws.Cells.Item(15,2).Value = ' https://english.news.cn/20250501/c982e5a78cd48a9bb5e7f219def71c3/c.html '
$ws.Cells.Item(15,3).Value = ' 2025-05-01 '

$ws.Cells.Item(16,1).Value = ' The Straits Times '
$ws.Cells.Item(16,2).Value = ' https://www.straitstimes.com/asia/east-asia/china-singapore-joint-ai-laboratory-opens-in-shenzhen '
$ws.Cells.Item(16,3).Value = ' 2025-04-30 '

$ws.Cells.Item(17,1).Value = ' China Daily '
$ws.Cells.Item(17,2).Value = ' https://www.chinadaily.com.cn/a/202504/27/WS66096a9a310dbb0113778c6.html '
$ws.Cells.Item(17,3).Value = ' 2025-04-27 '

$ws.Cells.Item(18,1).Value = ' South China Morning Post '
$ws.Cells.Item(18,2).Value = ' https://www.scmp.com/tech/big-tech/article/3254867/chinese-tech-giant-baidu-launches-156-qubit-quantum-computer-claiming-superior-performance '
$ws.Cells.Item(18,3).Value = ' 2025-04-28 '

$ws.Cells.Item(19,1).Value = ' Reuters Healthcare '
$ws.Cells.Item(19,2).Value = ' https://www.reuters.com/business/healthcare-pharmaceuticals/china-updates-biotech-regulatory-framework-accelerate-gene-therapy-approvals-2025-04-29/ '
$ws.Cells.Item(19,3).Value = ' 2025-04-29 '

$ws.Cells.Item(20,1).Value = ' Global Times '
$ws.Cells.Item(20,2).Value = ' https://www.globaltimes.cn/page/202504/1309825.shtml '
$ws.Cells.Item(20,3).Value = ' 2025-04-28 '

$ws.Cells.Item(21,1).Value = ' South China Morning Post '
$ws.Cells.Item(21,2).Value = ' https://www.scmp.com/news/china/science/article/3254832/china-announces-enhanced-talent-introduction-policy-expanded-visa-program '
$ws.Cells.Item(21,3).Value = ' 2025-04-27 '

$ws.Cells.Item(22,1).Value = ' Xinhua News Agency '
$ws.Cells.Item(22,2).Value = ' https://english.news.cn/20250429/a45fc6ee98f94d99b37d3b0cf3f2b64/c.html '
$ws.Cells.Item(22,3).Value = ' 2025-04-29 '

$ws.Rows("23:32").Delete()

# --- Executive Summary ---
$ws = $wb.Worksheets.Item("Executive Summary")
$execSummaryText = @'
Five Most Impactful News Summaries:

1. China's Innovation-Driven Development Strategy 2025-2035: The Chinese State Council has released a comprehensive ten-year innovation strategy that prioritizes technological self-reliance and indigenous innovation capabilities in critical sectors including semiconductors, advanced materials, and artificial intelligence. The strategy introduces a "2+X" framework focusing on core technologies and supporting ecosystem development, with targets to increase R&D spending to 3% of GDP by 2030. This represents China's most comprehensive science and technology roadmap to date and directly addresses current challenges in accessing certain foreign technologies.

2. National Artificial Intelligence Laboratory Network: China has established a national AI laboratory network comprising 15 specialized facilities focused on critical infrastructure protection and implementing the country's recently released ethical guidelines for AI. This initiative, falling under China's 14th Five-Year Plan for AI development, represents a significant organizational milestone in China's AI ecosystem, balancing innovation with security and ethical governance considerations. The network will coordinate research across multiple domains including healthcare, transportation, finance, and energy systems.

3. New Quality Productivity Manufacturing Initiative: The National Development and Reform Commission has announced a $150 billion initiative to establish 25 "New Quality Productivity" manufacturing hubs across China. This massive investment aims to upgrade China's industrial base with advanced digital manufacturing capabilities, smart factories, and integrated supply chains. Focusing initially on semiconductors, new energy vehicles, aerospace, and biotechnology, the program is expected to create 1.2 million high-skilled jobs and significantly advance China's position in high-value manufacturing.

4. Quantum Communication Network Expansion: Chinese researchers have achieved a major breakthrough by extending their secure quantum communication network to over 100 cities nationwide, demonstrating "quantum advantage" in secure communications using entangled photons. This represents the world's largest functional quantum communication infrastructure, with plans to extend services to Southeast Asia through the Digital Silk Road initiative, potentially revolutionizing secure communications and giving China a significant lead in quantum technologies.

5. Science and Technology Security Framework: China has introduced a comprehensive Science and Technology Security Framework with new export controls on AI chipsets and quantum technologies. The policy reflects growing concerns about technology containment strategies by Western nations while emphasizing the protection of critical indigenous innovations. This framework signals China's more assertive stance on technology sovereignty while maintaining its "fair and reciprocal" approach to international technology cooperation.
'@
$ws.Cells.Item(2,1).Value = $execSummaryText

# --- Cooccurrence ---
$ws = $wb.Worksheets.Item("Cooccurrence")
$ws.Rows("1:2").Delete()


# --- Associations ---
$ws = $wb.Worksheets.Item("Associations")
$ws.Cells.Item(2,1).Value = 'Innovation-driven development'
$ws.Cells.Item(2,2).Value = 1

$ws.Cells.Item(3,1).Value = 'New quality productivity'
$ws.Cells.Item(3,2).Value = 1

$ws.Cells.Item(4,1).Value = 'Quantum communication'
$ws.Cells.Item(4,2).Value = 1

$ws.Cells.Item(5,1).Value = 'Biotechnology'
$ws.Cells.Item(5,2).Value = 2

$ws.Cells.Item(6,1).Value = 'Semiconductor packaging'
$ws.Cells.Item(6,2).Value = 1

$ws.Cells.Item(7,1).Value = 'Science and technology innovation policy'
$ws.Cells.Item(7,2).Value = 1

$ws.Cells.Item(8,1).Value = 'Strategic emerging industries'
$ws.Cells.Item(8,2).Value = 1

$ws.Cells.Item(9,1).Value = 'Hydrogen energy storage'
$ws.Cells.Item(9,2).Value = 1

$ws.Cells.Item(10,1).Value = 'Science and technology system reform'
$ws.Cells.Item(10,2).Value = 1

$ws.Cells.Item(11,1).Value = 'Industry-university-research cooperation'
$ws.Cells.Item(11,2).Value = 1

$ws.Cells.Item(12,1).Value = 'Quantum computing cloud platform'
$ws.Cells.Item(12,2).Value = 1

$ws.Cells.Item(13,1).Value = 'High-end equipment manufacturing'
$ws.Cells.Item(13,2).Value = 1

$ws.Cells.Item(14,1).Value = 'Talent introduction policy'
$ws.Cells.Item(14,2).Value = 1

$ws.Cells.Item(15,1).Value = 'Science and technology cooperation agreement'
$ws.Cells.Item(15,2).Value = 1

